$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 111632005
$ws.Range("B2").Value = 89369
$ws.Range("D2").Value = 'LC'
$ws.Range("E2").Value = 5447
$ws.Range("F2").Value = 'Vedticka'
$ws.Range("G2").Value = 'Fuscoporia viticola'
$ws.Range("H2").Value = '(Schwein.) Murrill'
$ws.Range("Q2").Value = 679445.4687985049
$ws.Range("R2").Value = 6612665.387322281

# Row 3
$ws.Range("A3").Value = 111632616
$ws.Range("B3").Value = 89423
$ws.Range("D3").Value = 'NT'
$ws.Range("E3").Value = 5432
$ws.Range("F3").Value = 'Granticka'
$ws.Range("G3").Value = 'Porodaedalea chrysoloma'
$ws.Range("H3").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("Q3").Value = 679697.9675722983
$ws.Range("R3").Value = 6612862.479549611

# Row 4
$ws.Range("A4").Value = 111632126
$ws.Range("B4").Value = 90687
$ws.Range("D4").Value = 'LC'
$ws.Range("E4").Value = 5964
$ws.Range("F4").Value = 'Fjällig taggsvamp s.str.'
$ws.Range("G4").Value = 'Sarcodon imbricatus s.str.'
$ws.Range("H4").Value = '(L.:Fr.) P.Karst.'
$ws.Range("Q4").Value = 679400.8667491183
$ws.Range("R4").Value = 6612685.041705586

# Row 5
$ws.Range("A5").Value = 111632455
$ws.Range("B5").Value = 90687
$ws.Range("D5").Value = 'LC'
$ws.Range("E5").Value = 5964
$ws.Range("F5").Value = 'Fjällig taggsvamp s.str.'
$ws.Range("G5").Value = 'Sarcodon imbricatus s.str.'
$ws.Range("H5").Value = '(L.:Fr.) P.Karst.'
$ws.Range("Q5").Value = 679622.1256333978
$ws.Range("R5").Value = 6612739.744341305

# Row 6
$ws.Range("A6").Value = 111631648
$ws.Range("B6").Value = 90678
$ws.Range("D6").Value = 'LC'
$ws.Range("E6").Value = 4366
$ws.Range("F6").Value = 'Skarp dropptaggsvamp'
$ws.Range("G6").Value = 'Hydnellum peckii'
$ws.Range("H6").Value = 'Banker'
$ws.Range("Q6").Value = 679488.9265337941
$ws.Range("R6").Value = 6612786.06067825

# Row 7
$ws.Range("A7").Value = 111631780
$ws.Range("B7").Value = 90687
$ws.Range("D7").Value = 'LC'
$ws.Range("E7").Value = 5964
$ws.Range("F7").Value = 'Fjällig taggsvamp s.str.'
$ws.Range("G7").Value = 'Sarcodon imbricatus s.str.'
$ws.Range("H7").Value = '(L.:Fr.) P.Karst.'
$ws.Range("Q7").Value = 679488.9265337941
$ws.Range("R7").Value = 6612786.06067825

# Row 8
$ws.Range("A8").Value = 111631820
$ws.Range("B8").Value = 88819
$ws.Range("D8").Value = 'LC'
$ws.Range("E8").Value = 5685
$ws.Range("F8").Value = 'Gullgröppa'
$ws.Range("G8").Value = 'Pseudomerulius aureus'
$ws.Range("H8").Value = '(Fr.) Jülich'
$ws.Range("Q8").Value = 679445.4687985049
$ws.Range("R8").Value = 6612665.387322281

# Row 10
$ws.Range("A10").Value = 111631615
$ws.Range("B10").Value = 90687
$ws.Range("D10").Value = 'LC'
$ws.Range("E10").Value = 5964
$ws.Range("F10").Value = 'Fjällig taggsvamp s.str.'
$ws.Range("G10").Value = 'Sarcodon imbricatus s.str.'
$ws.Range("H10").Value = '(L.:Fr.) P.Karst.'
$ws.Range("Q10").Value = 679347.335090697
$ws.Range("R10").Value = 6612732.142182259

# Row 11
$ws.Range("A11").Value = 111631555
$ws.Range("B11").Value = 90658
$ws.Range("D11").Value = 'NT'
$ws.Range("E11").Value = 4361
$ws.Range("F11").Value = 'Orange taggsvamp'
$ws.Range("G11").Value = 'Hydnellum aurantiacum'
$ws.Range("H11").Value = '(Batsch:Fr.) P.Karst.'
$ws.Range("Q11").Value = 679347.335090697
$ws.Range("R11").Value = 6612732.142182259

# Row 12
$ws.Range("A12").Value = 111632936
$ws.Range("B12").Value = 90658
$ws.Range("D12").Value = 'NT'
$ws.Range("E12").Value = 4361
$ws.Range("F12").Value = 'Orange taggsvamp'
$ws.Range("G12").Value = 'Hydnellum aurantiacum'
$ws.Range("H12").Value = '(Batsch:Fr.) P.Karst.'
$ws.Range("Q12").Value = 679389.9201578975
$ws.Range("R12").Value = 6612881.656256998

# Row 13
$ws.Range("A13").Value = 111632295
$ws.Range("B13").Value = 89419
$ws.Range("D13").Value = 'NT'
$ws.Range("E13").Value = 1204
$ws.Range("F13").Value = 'Gränsticka'
$ws.Range("G13").Value = 'Phellopilus nigrolimitatus'
$ws.Range("H13").Value = '(Romell) Niemelä, T.Wagner & M.Fisch.'
$ws.Range("Q13").Value = 679615.9172154681
$ws.Range("R13").Value = 6612689.273485693
